$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.900.09'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '2.913.95'
$ws.Range("E3").Value = '  -0.38%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.505'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = '2.912.51'
$ws.Range("E9").Value = '  -0.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.01'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.24%  '

$ws.Range("E11").Value = '  +7.23%  '

$ws.Range("E12").Value = '  -1.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000238'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.51'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.48%  '

$ws.Range("E15").Value = '  -1.24%  '

$ws.Range("D16").Value = '3.400.10'
$ws.Range("E16").Value = '  -0.26%  '

$ws.Range("D17").Value = '61.915.25'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("D19").Value = '2.917.83'
$ws.Range("E19").Value = '  +0.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '434.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.659'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.36%  '

$ws.Range("E23").Value = '  -1.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.43%  '

$ws.Range("E26").Value = '  +1.61%  '

$ws.Range("E27").Value = '  -1.01%  '

$ws.Range("E28").Value = '  -0.06%  '

$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000106'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +24.14%  '

$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.64%  '

$ws.Range("E31").Value = '  -0.58%  '

$ws.Range("E32").Value = '  +0.75%  '

$ws.Range("E33").Value = '  +4.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.61%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.979'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.77%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.54'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.03%  '

$ws.Range("E40").Value = '  +3.02%  '

$ws.Range("E41").Value = '  -1.19%  '

$ws.Range("E42").Value = '  -1.53%  '

$ws.Range("E43").Value = '  +0.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.85%  '

$ws.Range("D45").Value = '2.699.70'
$ws.Range("E45").Value = '  +0.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '134.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.95%  '

$ws.Range("E47").Value = '  +1.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '347.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.24%  '

$ws.Range("E50").Value = '  +0.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.64%  '

